$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for United States in column F
$ws.Range("F1").Value = "United States"

# Update data values for columns B through F, rows 2-16
$ws.Range("B2").Value = 0.616858237547893
$ws.Range("C2").Value = 0.676258992805755
$ws.Range("D2").Value = 0.960199004975124
$ws.Range("E2").Value = 0.953307392996109
$ws.Range("F2").Value = 0.655973451327434
$ws.Range("B3").Value = 0.845559845559846
$ws.Range("C3").Value = 0.844594594594595
$ws.Range("D3").Value = 0.869158878504673
$ws.Range("E3").Value = 0.92057761732852
$ws.Range("F3").Value = 0.840262582056893
$ws.Range("B4").Value = 0.539748953974895
$ws.Range("C4").Value = 0.815181518151815
$ws.Range("D4").Value = 0.802690582959641
$ws.Range("E4").Value = 0.61003861003861
$ws.Range("F4").Value = 0.803474484256243
$ws.Range("B5").Value = 0.895752895752896
$ws.Range("C5").Value = 0.807817589576547
$ws.Range("D5").Value = 0.934579439252336
$ws.Range("E5").Value = 0.881294964028777
$ws.Range("F5").Value = 0.88628762541806
$ws.Range("B6").Value = 0.686507936507937
$ws.Range("C6").Value = 0.762345679012346
$ws.Range("D6").Value = 0.762931034482759
$ws.Range("E6").Value = 0.807547169811321
$ws.Range("F6").Value = 0.574786324786325
$ws.Range("B7").Value = 0.697247706422018
$ws.Range("C7").Value = 0.545454545454545
$ws.Range("D7").Value = 0.684210526315789
$ws.Range("E7").Value = 0.428044280442804
$ws.Range("F7").Value = 0.598941798941799
$ws.Range("B8").Value = 0.65843621399177
$ws.Range("C8").Value = 0.79672131147541
$ws.Range("D8").Value = 0.900473933649289
$ws.Range("E8").Value = 0.786290322580645
$ws.Range("F8").Value = 0.695966907962771
$ws.Range("B9").Value = 0.844444444444444
$ws.Range("C9").Value = 0.86084142394822
$ws.Range("D9").Value = 0.800947867298578
$ws.Range("E9").Value = 0.812
$ws.Range("F9").Value = 0.813928182807399
$ws.Range("B10").Value = 0.581395348837209
$ws.Range("C10").Value = 0.533546325878594
$ws.Range("D10").Value = 0.623853211009174
$ws.Range("E10").Value = 0.690839694656489
$ws.Range("F10").Value = 0.619895287958115
$ws.Range("B11").Value = 0.79746835443038
$ws.Range("C11").Value = 0.776073619631902
$ws.Range("D11").Value = 0.829596412556054
$ws.Range("E11").Value = 0.8
$ws.Range("F11").Value = 0.665938864628821
$ws.Range("B12").Value = 0.760504201680672
$ws.Range("C12").Value = 0.778156996587031
$ws.Range("D12").Value = 0.831683168316832
$ws.Range("E12").Value = 0.858921161825726
$ws.Range("F12").Value = 0.798553719008264
$ws.Range("B13").Value = 0.839449541284404
$ws.Range("C13").Value = 0.8561872909699
$ws.Range("D13").Value = 0.854077253218884
$ws.Range("E13").Value = 0.772549019607843
$ws.Range("F13").Value = 0.717127071823204
$ws.Range("B14").Value = 0.780392156862745
$ws.Range("C14").Value = 0.853820598006645
$ws.Range("D14").Value = 0.847290640394089
$ws.Range("E14").Value = 0.792307692307692
$ws.Range("F14").Value = 0.793926247288503
$ws.Range("B15").Value = 0.788235294117647
$ws.Range("C15").Value = 0.777003484320557
$ws.Range("D15").Value = 0.84549356223176
$ws.Range("E15").Value = 0.743295019157088
$ws.Range("F15").Value = 0.725690890481064
$ws.Range("B16").Value = 0.763565891472868
$ws.Range("C16").Value = 0.752380952380952
$ws.Range("D16").Value = 0.722943722943723
$ws.Range("E16").Value = 0.493877551020408
$ws.Range("F16").Value = 0.5917225950783
